$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "'307.59"
$ws.Range("E2").Value = "'-0.03%"
$ws.Range("D3").Value = "'41.07"
$ws.Range("E3").Value = "'0.04%"
$ws.Range("D4").Value = "'5.245"
$ws.Range("E4").Value = "'2.37%"
$ws.Range("D5").Value = "'0.07663"
$ws.Range("E5").Value = "'0.57%"
$ws.Range("D6").Value = "'1.624"
$ws.Range("E6").Value = "'0.31%"
$ws.Range("D7").Value = "'0.9187"
$ws.Range("E7").Value = "'1.89%"
$ws.Range("D8").Value = "'2.427"
$ws.Range("E8").Value = "'-3.21%"
$ws.Range("D9").Value = "'0.1226"
$ws.Range("E9").Value = "'12.44%"
$ws.Range("D10").Value = "'0.1824"
$ws.Range("E10").Value = "'2.71%"
$ws.Range("D11").Value = "'0.09081"
$ws.Range("E11").Value = "'-0.86%"
$ws.Range("D12").Value = "'0.04256"
$ws.Range("E12").Value = "'1.44%"
$ws.Range("E13").Value = "'-0.03%"
$ws.Range("D14").Value = "'0.001261"
$ws.Range("E14").Value = "'0.53%"
$ws.Range("D15").Value = "'0.005786"
$ws.Range("E15").Value = "'-0.12%"
$ws.Range("D17").Value = "'3.353"
$ws.Range("E17").Value = "'-0.11%"
$ws.Range("E18").Value = "'1.16%"
$ws.Range("E19").Value = "'1.22%"
$ws.Range("D20").Value = "'7.319"
$ws.Range("E20").Value = "'11.62%"
$ws.Range("D21").Value = "'0.1384"
$ws.Range("E21").Value = "'1.47%"
$ws.Range("E22").Value = "'2.87%"
$ws.Range("D23").Value = "'0.04075"
$ws.Range("E23").Value = "'0.05%"
$ws.Range("E24").Value = "'3.36%"
$ws.Range("D25").Value = "'0.004359"
$ws.Range("E25").Value = "'8.92%"
$ws.Range("E26").Value = "'-2.12%"
$ws.Range("D38").Value = "'0.02473"
$ws.Range("E38").Value = "'3.56%"
$ws.Range("D39").Value = "'0.05277"
$ws.Range("E39").Value = "'1.85%"
$ws.Range("D40").Value = "'0.007849"
$ws.Range("E40").Value = "'1.05%"
$ws.Range("D41").Value = "'0.1313"
$ws.Range("E41").Value = "'0.99%"
$ws.Range("D42").Value = "'0.006559"
$ws.Range("E42").Value = "'-5.81%"
$ws.Range("E43").Value = "'-1.86%"
$ws.Range("D44").Value = "'0.007668"
$ws.Range("E44").Value = "'-10.24%"
$ws.Range("D45").Value = "'0.3057"
$ws.Range("E45").Value = "'-0.50%"
$ws.Range("D46").Value = "'0.00006716"
$ws.Range("E46").Value = "'-2.48%"
$ws.Range("E47").Value = "'0.15%"
$ws.Range("D48").Value = "'0.4390"
$ws.Range("E48").Value = "'3,613.66%"
$ws.Range("E49").Value = "'-2.39%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'0.15%"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("E51").Value = "'0.15%"
